# The template used to ship with three blank tabs (Sheet1/Sheet2/Sheet3).
# Only Sheet1 actually holds the Community Learning report, so drop the
# two unused tabs and give the remaining one a descriptive name.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()

$wb.Worksheets.Item("Sheet1").Name = "Community Learning Report"

$ws = $wb.Worksheets.Item("Community Learning Report")

# Re-assert the print area so the workbook-level defined name
# ('_xlnm.Print_Area') is rewritten against the new sheet name instead of
# the stale "Sheet1" reference.
$ws.PageSetup.PrintArea = "`$A`$1:`$I`$42"
